$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row formatting: center (horizontal + vertical) A1:C1 ---
$headerRange = $ws.Range("A1:C1")
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4108     # xlCenter

# --- New data row for LeetCode 125 ---
$ws.Range("A2").Value = "Two Pointers"
$ws.Range("B2").Value = "125 - Valid Palindrome"
$ws.Range("C2").Value = "Python strings have several methods to check for their type of contents. Start a pointer at each side of the string and check for equality."

# Highlight the problem-name cell with the theme "Green, Accent 6" fill
$ws.Range("B2").Interior.ThemeColor = 10   # xlThemeColorAccent6 -> theme index 9

# --- Column widths to fit the new content ---
$ws.Columns.Item(1).ColumnWidth = 19.28515625
$ws.Columns.Item(2).ColumnWidth = 23.5703125
$ws.Columns.Item(3).ColumnWidth = 139.7109375

# --- Selection moves to C6 ---
$ws.Range("C6").Select()
